$d = $word.ActiveDocument

# Locate the paragraph that ends the bibliography entry ("...MacGrall-Hill") -
# this paragraph must be kept untouched.
$anchorIndex = -1
# Locate the paragraph containing the trailing "Powered by Jekyll" site-footer
# text - this paragraph (and everything between it and the anchor) is removed.
$footerIndex = -1

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*MacGrall-Hill*") {
        $anchorIndex = $i
    }
    if ($t -like "*Powered by Jekyll*") {
        $footerIndex = $i
    }
}

if ($anchorIndex -gt 0 -and $footerIndex -gt $anchorIndex) {
    $startPara = $d.Paragraphs.Item($anchorIndex + 1)
    $endPara = $d.Paragraphs.Item($footerIndex)
    $r = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $r.Delete()
}
